# Insert a new weekly price-report row for "Poroto verde" at row 111.
# This shifts the existing rows 111-131 down to rows 112-132 (unchanged),
# and populates the newly inserted row 111 with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 111, pushing old rows 111..131 to 112..132.
$ws.Rows("111:111").Insert()

# Fill in the new row 111 with the new weekly data point.
$ws.Cells.Item(111, 1).Value2 = 8
$ws.Cells.Item(111, 2).Value2 = "Terminal La Palmera de La Serena"
$ws.Cells.Item(111, 3).Value2 = "Coquimbo"
$ws.Cells.Item(111, 4).Value2 = 44505
$ws.Cells.Item(111, 5).Value2 = 4
$ws.Cells.Item(111, 6).Value2 = 100112031
$ws.Cells.Item(111, 7).Value2 = "Poroto verde"
$ws.Cells.Item(111, 8).Value2 = "Magnum"
$ws.Cells.Item(111, 9).Value2 = "Primera"
$ws.Cells.Item(111, 10).Value2 = 540
$ws.Cells.Item(111, 11).Value2 = 43000
$ws.Cells.Item(111, 12).Value2 = 44000
$ws.Cells.Item(111, 13).Value2 = 43500
$ws.Cells.Item(111, 14).Value2 = "`$/malla 25 kilos"
$ws.Cells.Item(111, 15).Value2 = "Provincia de Limarí"
$ws.Cells.Item(111, 16).Value2 = 1740
$ws.Cells.Item(111, 17).Value2 = 25
$ws.Cells.Item(111, 18).Value2 = "Hortaliza"
